$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value2 = 5065.4814
$ws.Range("I100").Value2 = 6098.421
$ws.Range("J100").Value2 = 2612.25
$ws.Range("K100").Value2 = 6098.421
$ws.Range("L100").Value2 = 2612.25
$ws.Range("M100").Value2 = -5557.421
$ws.Range("N100").Value2 = -3694.25

$ws.Range("H133").Value2 = 34000
$ws.Range("J133").Value2 = 34000
$ws.Range("L133").Value2 = 34000
$ws.Range("N133").Value2 = -44120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value2 = 7000
$ws.Range("J6").Value2 = 11000
$ws.Range("L6").Value2 = 11000
$ws.Range("N6").Value2 = -11346

$ws.Range("H74").Value2 = 17858436
$ws.Range("I74").Value2 = 19232050
$ws.Range("J74").Value2 = 1450
$ws.Range("K74").Value2 = 19232050
$ws.Range("L74").Value2 = 1450
$ws.Range("M74").Value2 = -19231176
$ws.Range("N74").Value2 = -3198

$ws.Range("H77").Value2 = 17858436
$ws.Range("I77").Value2 = 19232050
$ws.Range("J77").Value2 = 1450
$ws.Range("K77").Value2 = 96160250
$ws.Range("L77").Value2 = 7250
$ws.Range("M77").Value2 = -96155882
$ws.Range("N77").Value2 = -15986

$ws.Range("H132").Value2 = 32263024
$ws.Range("I132").Value2 = 45455908
$ws.Range("J132").Value2 = 13755.333
$ws.Range("K132").Value2 = 136367724
$ws.Range("L132").Value2 = 41265.999
$ws.Range("M132").Value2 = -136365194
$ws.Range("N132").Value2 = -46325.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 2140153.8
$ws.Range("I134").Value2 = 1078.5294
$ws.Range("K134").Value2 = 3235.5882
$ws.Range("M134").Value2 = -700.5881999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value2 = 637.4
$ws.Range("J10").Value2 = 702.6667
$ws.Range("L10").Value2 = 702.6667
$ws.Range("N10").Value2 = -980.6667

$ws.Range("H11").Value2 = 4006
$ws.Range("I11").Value2 = 0
$ws.Range("J11").Value2 = 4006
$ws.Range("K11").Value2 = 0
$ws.Range("L11").Value2 = 4006
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value2 = -4286

$ws.Range("H13").Value2 = 4005
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 4005
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 4005
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value2 = -4283

$ws.Range("H14").Value2 = 7777.375
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 7777.375
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 7777.375
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value2 = -8117.375

$ws.Range("H19").Value2 = 383.16666
$ws.Range("I19").Value2 = 383.16666
$ws.Range("K19").Value2 = 383.16666
$ws.Range("M19").Value2 = -213.16666

$ws.Range("H24").Value2 = 383.16666
$ws.Range("I24").Value2 = 383.16666
$ws.Range("K24").Value2 = 383.16666
$ws.Range("M24").Value2 = -213.16666

$ws.Range("H31").Value2 = 1585.091
$ws.Range("I31").Value2 = 1083.0667
$ws.Range("J31").Value2 = 2291.0625
$ws.Range("K31").Value2 = 1083.0667
$ws.Range("L31").Value2 = 2291.0625
$ws.Range("M31").Value2 = -788.0667000000001
$ws.Range("N31").Value2 = -2881.0625

$ws.Range("H34").Value2 = 1585.091
$ws.Range("I34").Value2 = 1083.0667
$ws.Range("J34").Value2 = 2291.0625
$ws.Range("K34").Value2 = 1083.0667
$ws.Range("L34").Value2 = 2291.0625
$ws.Range("M34").Value2 = -881.0667000000001
$ws.Range("N34").Value2 = -2695.0625

$ws.Range("H94").Value2 = 2651.6
$ws.Range("I94").Value2 = 1200.6666
$ws.Range("J94").Value2 = 3273.4285
$ws.Range("K94").Value2 = 1200.6666
$ws.Range("L94").Value2 = 3273.4285
$ws.Range("M94").Value2 = -749.6666
$ws.Range("N94").Value2 = -4175.4285

$ws.Range("H134").Value2 = 1301.2972
$ws.Range("I134").Value2 = 1364.96
$ws.Range("K134").Value2 = 4094.88
$ws.Range("M134").Value2 = -1559.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 55571.668
$ws.Range("I12").Value2 = 8
$ws.Range("J12").Value2 = 76942.30499999999
$ws.Range("K12").Value2 = 24
$ws.Range("L12").Value2 = 230826.915
$ws.Range("M12").Value2 = 149
$ws.Range("N12").Value2 = -231172.915

$ws.Range("H131").Value2 = 747.37
$ws.Range("I131").Value2 = 433.52942
$ws.Range("J131").Value2 = 811.6506000000001
$ws.Range("K131").Value2 = 1300.58826
$ws.Range("L131").Value2 = 2434.9518
$ws.Range("M131").Value2 = 3739.41174
$ws.Range("N131").Value2 = -12514.9518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 2424.1177
$ws.Range("I113").Value2 = 2471.4285
$ws.Range("J113").Value2 = 2391
$ws.Range("K113").Value2 = 2471.4285
$ws.Range("L113").Value2 = 2391
$ws.Range("M113").Value2 = -301.4285
$ws.Range("N113").Value2 = -6731

$ws.Range("H122").Value2 = 18524906
$ws.Range("I122").Value2 = 20839770
$ws.Range("K122").Value2 = 62519310
$ws.Range("M122").Value2 = -62516860

$ws.Range("H123").Value2 = 19313.133
$ws.Range("J123").Value2 = 19313.133
$ws.Range("L123").Value2 = 19313.133
$ws.Range("N123").Value2 = -24213.133

$ws.Range("H132").Value2 = 24999.6
$ws.Range("I132").Value2 = 2000
$ws.Range("J132").Value2 = 30749.5
$ws.Range("K132").Value2 = 6000
$ws.Range("L132").Value2 = 92248.5
$ws.Range("M132").Value2 = -3470
$ws.Range("N132").Value2 = -97308.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 2977521.8
$ws.Range("I46").Value2 = 4630225
$ws.Range("J46").Value2 = 2655.8
$ws.Range("K46").Value2 = 4630225
$ws.Range("L46").Value2 = 2655.8
$ws.Range("M46").Value2 = -4630037
$ws.Range("N46").Value2 = -3031.8

$ws.Range("H61").Value2 = 2383.1304
$ws.Range("I61").Value2 = 2312.625
$ws.Range("J61").Value2 = 2544.2856
$ws.Range("K61").Value2 = 2312.625
$ws.Range("L61").Value2 = 2544.2856
$ws.Range("M61").Value2 = -2110.625
$ws.Range("N61").Value2 = -2948.2856

$ws.Range("H88").Value2 = 10000
$ws.Range("I88").Value2 = 10000
$ws.Range("J88").Value2 = 0
$ws.Range("K88").Value2 = 10000
$ws.Range("L88").Value2 = 0
$ws.Range("M88").Value2 = -9572
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value2 = 10000
$ws.Range("I91").Value2 = 10000
$ws.Range("J91").Value2 = 0
$ws.Range("K91").Value2 = 10000
$ws.Range("L91").Value2 = 0
$ws.Range("M91").Value2 = -8518
$ws.Range("N91").ClearContents()

$ws.Range("H113").Value2 = 2383.1304
$ws.Range("I113").Value2 = 2312.625
$ws.Range("J113").Value2 = 2544.2856
$ws.Range("K113").Value2 = 2312.625
$ws.Range("L113").Value2 = 2544.2856
$ws.Range("M113").Value2 = -142.625
$ws.Range("N113").Value2 = -6884.2856
